# Upgrade all packages and take latest June 2024 data
# Update the "2024" section (rows 191-199) of Sheet1 with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 191: Richard (rank 1) - values updated, person unchanged
$ws.Cells.Item(191, 4).Value  = 35        # D
$ws.Cells.Item(191, 6).Value  = 35        # F
$ws.Cells.Item(191, 7).Value  = 108650    # G
$ws.Cells.Item(191, 8).Value  = 140       # H
$ws.Cells.Item(191, 9).Value  = 90        # I

# Row 192: Mark (rank 2)
$ws.Cells.Item(192, 2).Value  = "Mark"    # B
$ws.Cells.Item(192, 4).Value  = 25        # D
$ws.Cells.Item(192, 6).Value  = 25        # F
$ws.Cells.Item(192, 7).Value  = 74950     # G
$ws.Cells.Item(192, 8).Value  = 80        # H
$ws.Cells.Item(192, 9).Value  = 20        # I
$ws.Cells.Item(192, 11).Value = 361       # K

# Row 193: Andy (rank 3)
$ws.Cells.Item(193, 2).Value  = "Andy"    # B
$ws.Cells.Item(193, 4).Value  = 23        # D
$ws.Cells.Item(193, 6).Value  = 23        # F
$ws.Cells.Item(193, 7).Value  = 73650     # G
$ws.Cells.Item(193, 8).Value  = 60        # H
$ws.Cells.Item(193, 9).Value  = 10        # I
$ws.Cells.Item(193, 11).Value = 349       # K

# Row 194: Anthony (rank 4)
$ws.Cells.Item(194, 2).Value  = "Anthony" # B
$ws.Cells.Item(194, 4).Value  = 23        # D
$ws.Cells.Item(194, 6).Value  = 23        # F
$ws.Cells.Item(194, 7).Value  = 71750     # G
$ws.Cells.Item(194, 8).Value  = 40        # H
$ws.Cells.Item(194, 9).Value  = -10       # I
$ws.Cells.Item(194, 11).Value = 350       # K

# Row 195: Jon (rank 5)
$ws.Cells.Item(195, 2).Value  = "Jon"     # B
$ws.Cells.Item(195, 4).Value  = 18        # D
$ws.Cells.Item(195, 6).Value  = 18        # F
$ws.Cells.Item(195, 7).Value  = 50650     # G
$ws.Cells.Item(195, 8).Value  = 20        # H
$ws.Cells.Item(195, 9).Value  = -40       # I
$ws.Cells.Item(195, 11).Value = 357       # K

# Row 196: Pepe (rank 6)
$ws.Cells.Item(196, 2).Value  = "Pepe"    # B
$ws.Cells.Item(196, 4).Value  = 18        # D
$ws.Cells.Item(196, 6).Value  = 18        # F
$ws.Cells.Item(196, 7).Value  = 50100     # G
$ws.Cells.Item(196, 8).Value  = 40        # H
$ws.Cells.Item(196, 9).Value  = -10       # I
$ws.Cells.Item(196, 11).Value = 364       # K

# Row 197: Prashant (rank 7)
$ws.Cells.Item(197, 2).Value  = "Prashant" # B
$ws.Cells.Item(197, 4).Value  = 14         # D
$ws.Cells.Item(197, 6).Value  = 14         # F
$ws.Cells.Item(197, 7).Value  = 48250      # G
$ws.Cells.Item(197, 8).Value  = 60         # H
$ws.Cells.Item(197, 9).Value  = 20         # I
$ws.Cells.Item(197, 11).Value = 365        # K

# Row 198: Matt (rank 8)
$ws.Cells.Item(198, 2).Value  = "Matt"    # B
$ws.Cells.Item(198, 4).Value  = 12        # D
$ws.Cells.Item(198, 6).Value  = 12        # F
$ws.Cells.Item(198, 7).Value  = 55450     # G
$ws.Cells.Item(198, 8).Value  = 10        # H
$ws.Cells.Item(198, 9).Value  = -50       # I
$ws.Cells.Item(198, 11).Value = 362       # K

# Row 199: Maisy (rank 9)
$ws.Cells.Item(199, 2).Value  = "Maisy"   # B
$ws.Cells.Item(199, 4).Value  = 8         # D
$ws.Cells.Item(199, 6).Value  = 8         # F
$ws.Cells.Item(199, 7).Value  = 35700     # G
$ws.Cells.Item(199, 8).Value  = 20        # H
$ws.Cells.Item(199, 9).Value  = -30       # I
$ws.Cells.Item(199, 11).Value = 360       # K
